$d = $word.ActiveDocument

# 1. Update the "Generated" timestamp.
$d.Content.Find.Execute(
    "Generated: 2026-02-07 12:11:40", $true, $false, $false, $false, $false,
    $true, 1, $false, "Generated: 2026-02-10 15:36:36", 2) | Out-Null

# 2. Update the summary table values.
#    Row order: Total Paragraphs / Auto-Applied / Needs Review / Auto-Apply Rate
#    (Note: once a table is touched, Document.Paragraphs indexing becomes
#    unreliable in this runtime, so every later step below is driven off
#    Range objects - e.g. $d.Content - rather than the Paragraphs collection.)
$t = $d.Tables.Item(1)
$t.Cell(2, 2).Range.Text = "445"
$t.Cell(3, 2).Range.Text = "1"
$t.Cell(4, 2).Range.Text = "99.8%"

# 3. Swap the "all clear" heading/paragraph for the "needs review" versions.
$d.Content.Find.Execute(
    "All Items Auto-Applied", $true, $false, $false, $false, $false,
    $true, 1, $false, "Items Requiring Review", 2) | Out-Null

$noReviewText = "All paragraphs were classified with high confidence (" + [char]8805 + "85%). No manual review required."
$d.Content.Find.Execute(
    $noReviewText, $true, $false, $false, $false, $false,
    $true, 1, $false, "The following 1 items have confidence below 85% and require human review.", 2) | Out-Null

# 4. Append the detail block for the single item that needs review.
$rng = $d.Content
$rng.Collapse(0)

$detailText = "`r`r" + "Paragraph 169`r" + `
    'Text: "<DFIG2.1>Figure 2.1: Anatomy of a CAR Construct. Adapted from: Acharya, U.H.; Walter, R.B. Chimeric ..."' + "`r" + `
    "Suggested Tag: REF-U (Confidence: 80%)`r" + `
    "Reasoning: Figure legend with source attribution`r"

$rng.InsertAfter($detailText)

# Bold just the "Paragraph 169" heading line.
$boldRng = $d.Content
$boldRng.Find.Execute("Paragraph 169", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boldRng.Bold = 1
